# Add data for 2024-04-13
# Updates year-to-date (2024, column K) violent-crime counts for the newly
# reported day across the citywide summary, the by-neighborhood rollup, and
# every affected neighborhood detail sheet. A couple of sheets also carry a
# small correction to the prior-year (2023, column J) figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1996
$ws.Range("K3").Value = 1928
$ws.Range("J4").Value = 1802
$ws.Range("K4").Value = 407
$ws.Range("K5").Value = 126
$ws.Range("K6").Value = 2453
$ws.Range("J7").Value = 29273
$ws.Range("K7").Value = 6910

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 202
$ws.Range("K8").Value = 463
$ws.Range("K9").Value = 28
$ws.Range("K15").Value = 63
$ws.Range("K19").Value = 191
$ws.Range("K20").Value = 149
$ws.Range("K24").Value = 22
$ws.Range("K27").Value = 80
$ws.Range("J29").Value = 1557
$ws.Range("K29").Value = 340
$ws.Range("K33").Value = 281
$ws.Range("K34").Value = 42
$ws.Range("K36").Value = 82
$ws.Range("K37").Value = 230
$ws.Range("K40").Value = 14
$ws.Range("K41").Value = 63
$ws.Range("K42").Value = 240
$ws.Range("K44").Value = 67
$ws.Range("K48").Value = 82
$ws.Range("K49").Value = 41
$ws.Range("K50").Value = 39
$ws.Range("K52").Value = 187
$ws.Range("K54").Value = 123
$ws.Range("K55").Value = 71
$ws.Range("K60").Value = 48
$ws.Range("J63").Value = 94
$ws.Range("K65").Value = 168
$ws.Range("K67").Value = 266
$ws.Range("K68").Value = 18
$ws.Range("K71").Value = 18
$ws.Range("K72").Value = 32
$ws.Range("K73").Value = 68
$ws.Range("K76").Value = 101
$ws.Range("K77").Value = 49
$ws.Range("K78").Value = 89
$ws.Range("K79").Value = 181
$ws.Range("J83").Value = 592
$ws.Range("K83").Value = 149
$ws.Range("K84").Value = 47
$ws.Range("J85").Value = 1195
$ws.Range("K85").Value = 345
$ws.Range("K86").Value = 46
$ws.Range("K88").Value = 91
$ws.Range("K91").Value = 66
$ws.Range("K92").Value = 31
$ws.Range("K96").Value = 97
$ws.Range("K97").Value = 62
$ws.Range("K99").Value = 124
$ws.Range("J101").Value = 29273
$ws.Range("K101").Value = 6910

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 16
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 70
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 433
$ws.Range("K6").Value = 83
$ws.Range("J7").Value = 1195
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 133
$ws.Range("K3").Value = 133
$ws.Range("K6").Value = 162
$ws.Range("K7").Value = 463

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 60
$ws.Range("J4").Value = 24
$ws.Range("K6").Value = 32
$ws.Range("J7").Value = 592
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 77
$ws.Range("K3").Value = 108
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 47
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 266

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 88
$ws.Range("K3").Value = 112
$ws.Range("J4").Value = 84
$ws.Range("K4").Value = 17
$ws.Range("J7").Value = 1557
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 60
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 21
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 15
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 14
